$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New columns: K = "Typist", L = "Typist QC" with matching data rows,
# mirroring the existing header/client columns already on the sheet
# (new client "Accurate" onboarding columns).
$ws.Range("K1").Value = "Typist"
$ws.Range("L1").Value = "Typist QC"
$ws.Range("K2").Value = "SIPL0102"
$ws.Range("L2").Value = "SIPL0103"
$ws.Range("K3").Value = "SIPL0102"
$ws.Range("L3").Value = "SIPL0103"

# Match formatting of the existing header row / data rows for the new columns.
$ws.Range("J1").Copy()
$ws.Range("K1:L1").PasteSpecial(-4122)

$ws.Range("J2").Copy()
$ws.Range("K2:L2").PasteSpecial(-4122)
$ws.Range("K3:L3").PasteSpecial(-4122)

# Re-apply values (PasteSpecial of formats only shouldn't touch them, but
# make sure the correct text/shared-strings are set regardless of paste mode).
$ws.Range("K1").Value = "Typist"
$ws.Range("L1").Value = "Typist QC"
$ws.Range("K2").Value = "SIPL0102"
$ws.Range("L2").Value = "SIPL0103"
$ws.Range("K3").Value = "SIPL0102"
$ws.Range("L3").Value = "SIPL0103"

# Restore the active selection to match the saved view state.
$ws.Range("L12").Select()
